$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of (row, col) -> new text, in document order matching the diff.
# The populated rows in this table are 1, 5, 9, 13, 17 (5 cells each).
$updates = @(
    @{ Row = 1;  Col = 1; Text = "90÷8=11, 2" },
    @{ Row = 1;  Col = 2; Text = "40÷5=8, 0" },
    @{ Row = 1;  Col = 3; Text = "38÷8=4, 6" },
    @{ Row = 1;  Col = 4; Text = "93÷3=31, 0" },
    @{ Row = 1;  Col = 5; Text = "80÷4=20, 0" },

    @{ Row = 5;  Col = 1; Text = "67÷9=7, 4" },
    @{ Row = 5;  Col = 2; Text = "50÷8=6, 2" },
    @{ Row = 5;  Col = 3; Text = "62÷4=15, 2" },
    @{ Row = 5;  Col = 4; Text = "48÷6=8, 0" },
    @{ Row = 5;  Col = 5; Text = "49÷9=5, 4" },

    @{ Row = 9;  Col = 1; Text = "74÷3=24, 2" },
    @{ Row = 9;  Col = 2; Text = "33÷8=4, 1" },
    @{ Row = 9;  Col = 3; Text = "35÷2=17, 1" },
    @{ Row = 9;  Col = 4; Text = "77÷8=9, 5" },
    @{ Row = 9;  Col = 5; Text = "44÷5=8, 4" },

    @{ Row = 13; Col = 1; Text = "37÷9=4, 1" },
    @{ Row = 13; Col = 2; Text = "75÷8=9, 3" },
    @{ Row = 13; Col = 3; Text = "84÷5=16, 4" },
    @{ Row = 13; Col = 4; Text = "32÷3=10, 2" },
    @{ Row = 13; Col = 5; Text = "80÷3=26, 2" },

    @{ Row = 17; Col = 1; Text = "21÷8=2, 5" },
    @{ Row = 17; Col = 2; Text = "15÷3=5, 0" },
    @{ Row = 17; Col = 3; Text = "86÷6=14, 2" },
    @{ Row = 17; Col = 4; Text = "59÷5=11, 4" },
    @{ Row = 17; Col = 5; Text = "63÷6=10, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}

Write-Output "done"
